$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"3.711779"
$ws.Range("H2").Value = [double]"11.135337"
$ws.Range("I2").Value = [double]"0.008539583513749102"
$ws.Range("J2").Value = [double]"0.008539583513749104"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.926994"
$ws.Range("N2").Value = [double]"2.780982"
$ws.Range("O2").Value = [double]"0.002566096653125693"
$ws.Range("P2").Value = [double]"0.002566096653125693"
$ws.Range("Q2").Value = [double]"3.440796862326"
$ws.Range("R2").Value = [double]"30.967171760934"
$ws.Range("S2").Value = [double]"2.191339667371892E-05"
$ws.Range("T2").Value = [double]"2.191339667371892E-05"

$ws.Range("G3").Value = [double]"3.711779"
$ws.Range("H3").Value = [double]"11.135337"
$ws.Range("I3").Value = [double]"0.008539583513749102"
$ws.Range("J3").Value = [double]"0.008539583513749104"
$ws.Range("M3").Value = [double]"93.12610233333334"
$ws.Range("N3").Value = [double]"279.378307"
$ws.Range("O3").Value = [double]"0.2577908589658698"
$ws.Range("P3").Value = [double]"0.2577908589658698"
$ws.Range("Q3").Value = [double]"345.6635109927176"
$ws.Range("R3").Value = [double]"3110.971598934459"
$ws.Range("S3").Value = [double]"0.002201426569220162"
$ws.Range("T3").Value = [double]"0.002201426569220162"

$ws.Range("G4").Value = [double]"3.711779"
$ws.Range("H4").Value = [double]"11.135337"
$ws.Range("I4").Value = [double]"0.008539583513749102"
$ws.Range("J4").Value = [double]"0.008539583513749104"
$ws.Range("M4").Value = [double]"264.9957936666667"
$ws.Range("N4").Value = [double]"794.9873809999999"
$ws.Range("O4").Value = [double]"0.733559029746061"
$ws.Range("P4").Value = [double]"0.733559029746061"
$ws.Range("Q4").Value = [double]"983.6058220202663"
$ws.Range("R4").Value = [double]"8852.452398182397"
$ws.Range("S4").Value = [double]"0.00626428859678125"
$ws.Range("T4").Value = [double]"0.006264288596781251"

$ws.Range("G5").Value = [double]"3.711779"
$ws.Range("H5").Value = [double]"11.135337"
$ws.Range("I5").Value = [double]"0.008539583513749102"
$ws.Range("J5").Value = [double]"0.008539583513749104"
$ws.Range("M5").Value = [double]"2.197830333333333"
$ws.Range("N5").Value = [double]"6.593490999999999"
$ws.Range("O5").Value = [double]"0.006084014634943477"
$ws.Range("P5").Value = [double]"0.006084014634943477"
$ws.Range("Q5").Value = [double]"8.157860476829665"
$ws.Range("R5").Value = [double]"73.42074429146699"
$ws.Range("S5").Value = [double]"5.195495107397158E-05"
$ws.Range("T5").Value = [double]"5.195495107397159E-05"

$ws.Range("I6").Value = [double]"0.00304268690962334"
$ws.Range("J6").Value = [double]"0.003042686909623341"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.926994"
$ws.Range("N6").Value = [double]"2.780982"
$ws.Range("O6").Value = [double]"0.002566096653125693"
$ws.Range("P6").Value = [double]"0.002566096653125693"
$ws.Range("Q6").Value = [double]"1.225969340872"
$ws.Range("R6").Value = [double]"11.033724067848"
$ws.Range("S6").Value = [double]"7.807828695293812E-06"
$ws.Range("T6").Value = [double]"7.807828695293815E-06"

$ws.Range("I7").Value = [double]"0.00304268690962334"
$ws.Range("J7").Value = [double]"0.003042686909623341"
$ws.Range("M7").Value = [double]"93.12610233333334"
$ws.Range("N7").Value = [double]"279.378307"
$ws.Range("O7").Value = [double]"0.2577908589658698"
$ws.Range("P7").Value = [double]"0.2577908589658698"
$ws.Range("Q7").Value = [double]"123.1612570260164"
$ws.Range("R7").Value = [double]"1108.451313234148"
$ws.Range("S7").Value = [double]"0.0007843768719960088"
$ws.Range("T7").Value = [double]"0.000784376871996009"

$ws.Range("I8").Value = [double]"0.00304268690962334"
$ws.Range("J8").Value = [double]"0.003042686909623341"
$ws.Range("M8").Value = [double]"264.9957936666667"
$ws.Range("N8").Value = [double]"794.9873809999999"
$ws.Range("O8").Value = [double]"0.733559029746061"
$ws.Range("P8").Value = [double]"0.733559029746061"
$ws.Range("Q8").Value = [double]"350.4625903677648"
$ws.Range("R8").Value = [double]"3154.163313309884"
$ws.Range("S8").Value = [double]"0.002231990457244338"
$ws.Range("T8").Value = [double]"0.002231990457244339"

$ws.Range("I9").Value = [double]"0.00304268690962334"
$ws.Range("J9").Value = [double]"0.003042686909623341"
$ws.Range("M9").Value = [double]"2.197830333333333"
$ws.Range("N9").Value = [double]"6.593490999999999"
$ws.Range("O9").Value = [double]"0.006084014634943477"
$ws.Range("P9").Value = [double]"0.006084014634943477"
$ws.Range("Q9").Value = [double]"2.906677502880444"
$ws.Range("R9").Value = [double]"26.160097525924"
$ws.Range("S9").Value = [double]"1.851175168769934E-05"
$ws.Range("T9").Value = [double]"1.851175168769935E-05"

$ws.Range("G10").Value = [double]"81.93664033333333"
$ws.Range("H10").Value = [double]"245.809921"
$ws.Range("I10").Value = [double]"0.1885092789636783"
$ws.Range("J10").Value = [double]"0.1885092789636784"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"0.6666666666666666"
$ws.Range("M10").Value = [double]"0.926994"
$ws.Range("N10").Value = [double]"2.780982"
$ws.Range("O10").Value = [double]"0.002566096653125693"
$ws.Range("P10").Value = [double]"0.002566096653125693"
$ws.Range("Q10").Value = [double]"75.954773969158"
$ws.Range("R10").Value = [double]"683.5929657224219"
$ws.Range("S10").Value = [double]"0.0004837330298318326"
$ws.Range("T10").Value = [double]"0.0004837330298318328"

$ws.Range("G11").Value = [double]"81.93664033333333"
$ws.Range("H11").Value = [double]"245.809921"
$ws.Range("I11").Value = [double]"0.1885092789636783"
$ws.Range("J11").Value = [double]"0.1885092789636784"
$ws.Range("M11").Value = [double]"93.12610233333334"
$ws.Range("N11").Value = [double]"279.378307"
$ws.Range("O11").Value = [double]"0.2577908589658698"
$ws.Range("P11").Value = [double]"0.2577908589658698"
$ws.Range("Q11").Value = [double]"7630.439952531528"
$ws.Range("R11").Value = [double]"68673.95957278374"
$ws.Range("S11").Value = [double]"0.04859596894708341"
$ws.Range("T11").Value = [double]"0.04859596894708342"

$ws.Range("G12").Value = [double]"81.93664033333333"
$ws.Range("H12").Value = [double]"245.809921"
$ws.Range("I12").Value = [double]"0.1885092789636783"
$ws.Range("J12").Value = [double]"0.1885092789636784"
$ws.Range("M12").Value = [double]"264.9957936666667"
$ws.Range("N12").Value = [double]"794.9873809999999"
$ws.Range("O12").Value = [double]"0.733559029746061"
$ws.Range("P12").Value = [double]"0.733559029746061"
$ws.Range("Q12").Value = [double]"21712.86503551188"
$ws.Range("R12").Value = [double]"195415.7853196069"
$ws.Range("S12").Value = [double]"0.1382826837747254"
$ws.Range("T12").Value = [double]"0.1382826837747255"

$ws.Range("G13").Value = [double]"81.93664033333333"
$ws.Range("H13").Value = [double]"245.809921"
$ws.Range("I13").Value = [double]"0.1885092789636783"
$ws.Range("J13").Value = [double]"0.1885092789636784"
$ws.Range("M13").Value = [double]"2.197830333333333"
$ws.Range("N13").Value = [double]"6.593490999999999"
$ws.Range("O13").Value = [double]"0.006084014634943477"
$ws.Range("P13").Value = [double]"0.006084014634943477"
$ws.Range("Q13").Value = [double]"180.0828335360234"
$ws.Range("R13").Value = [double]"1620.745501824211"
$ws.Range("S13").Value = [double]"0.001146893212037661"
$ws.Range("T13").Value = [double]"0.001146893212037662"

$ws.Range("G14").Value = [double]"0.4187683333333334"
$ws.Range("H14").Value = [double]"1.256305"
$ws.Range("I14").Value = [double]"0.0009634482967368267"
$ws.Range("J14").Value = [double]"0.0009634482967368269"
$ws.Range("K14").Value = [double]"2"
$ws.Range("L14").Value = [double]"0.6666666666666666"
$ws.Range("M14").Value = [double]"0.926994"
$ws.Range("N14").Value = [double]"2.780982"
$ws.Range("O14").Value = [double]"0.002566096653125693"
$ws.Range("P14").Value = [double]"0.002566096653125693"
$ws.Range("Q14").Value = [double]"0.38819573239"
$ws.Range("R14").Value = [double]"3.49376159151"
$ws.Range("S14").Value = [double]"2.472301449716021E-06"
$ws.Range("T14").Value = [double]"2.472301449716021E-06"

$ws.Range("G15").Value = [double]"0.4187683333333334"
$ws.Range("H15").Value = [double]"1.256305"
$ws.Range("I15").Value = [double]"0.0009634482967368267"
$ws.Range("J15").Value = [double]"0.0009634482967368269"
$ws.Range("M15").Value = [double]"93.12610233333334"
$ws.Range("N15").Value = [double]"279.378307"
$ws.Range("O15").Value = [double]"0.2577908589658698"
$ws.Range("P15").Value = [double]"0.2577908589658698"
$ws.Range("Q15").Value = [double]"38.99826266395944"
$ws.Range("R15").Value = [double]"350.984363975635"
$ws.Range("S15").Value = [double]"0.0002483681639849907"
$ws.Range("T15").Value = [double]"0.0002483681639849908"

$ws.Range("G16").Value = [double]"0.4187683333333334"
$ws.Range("H16").Value = [double]"1.256305"
$ws.Range("I16").Value = [double]"0.0009634482967368267"
$ws.Range("J16").Value = [double]"0.0009634482967368269"
$ws.Range("M16").Value = [double]"264.9957936666667"
$ws.Range("N16").Value = [double]"794.9873809999999"
$ws.Range("O16").Value = [double]"0.733559029746061"
$ws.Range("P16").Value = [double]"0.733559029746061"
$ws.Range("Q16").Value = [double]"110.9718468541339"
$ws.Range("R16").Value = [double]"998.7466216872049"
$ws.Range("S16").Value = [double]"0.0007067461977647617"
$ws.Range("T16").Value = [double]"0.0007067461977647619"

$ws.Range("G17").Value = [double]"0.4187683333333334"
$ws.Range("H17").Value = [double]"1.256305"
$ws.Range("I17").Value = [double]"0.0009634482967368267"
$ws.Range("J17").Value = [double]"0.0009634482967368269"
$ws.Range("M17").Value = [double]"2.197830333333333"
$ws.Range("N17").Value = [double]"6.593490999999999"
$ws.Range("O17").Value = [double]"0.006084014634943477"
$ws.Range("P17").Value = [double]"0.006084014634943477"
$ws.Range("Q17").Value = [double]"0.9203817456394444"
$ws.Range("R17").Value = [double]"8.283435710754999"
$ws.Range("S17").Value = [double]"5.861633537358219E-06"
$ws.Range("T17").Value = [double]"5.86163353735822E-06"

$ws.Range("G18").Value = [double]"276.259491"
$ws.Range("H18").Value = [double]"828.7784730000001"
$ws.Range("I18").Value = [double]"0.6355822894790661"
$ws.Range("J18").Value = [double]"0.6355822894790663"
$ws.Range("K18").Value = [double]"2"
$ws.Range("L18").Value = [double]"0.6666666666666666"
$ws.Range("M18").Value = [double]"0.926994"
$ws.Range("N18").Value = [double]"2.780982"
$ws.Range("O18").Value = [double]"0.002566096653125693"
$ws.Range("P18").Value = [double]"0.002566096653125693"
$ws.Range("Q18").Value = [double]"256.090890600054"
$ws.Range("R18").Value = [double]"2304.818015400486"
$ws.Range("S18").Value = [double]"0.001630965585818197"
$ws.Range("T18").Value = [double]"0.001630965585818197"

$ws.Range("G19").Value = [double]"276.259491"
$ws.Range("H19").Value = [double]"828.7784730000001"
$ws.Range("I19").Value = [double]"0.6355822894790661"
$ws.Range("J19").Value = [double]"0.6355822894790663"
$ws.Range("M19").Value = [double]"93.12610233333334"
$ws.Range("N19").Value = [double]"279.378307"
$ws.Range("O19").Value = [double]"0.2577908589658698"
$ws.Range("P19").Value = [double]"0.2577908589658698"
$ws.Range("Q19").Value = [double]"25726.96962942058"
$ws.Range("R19").Value = [double]"231542.7266647852"
$ws.Range("S19").Value = [double]"0.1638473043483026"
$ws.Range("T19").Value = [double]"0.1638473043483026"

$ws.Range("G20").Value = [double]"276.259491"
$ws.Range("H20").Value = [double]"828.7784730000001"
$ws.Range("I20").Value = [double]"0.6355822894790661"
$ws.Range("J20").Value = [double]"0.6355822894790663"
$ws.Range("M20").Value = [double]"264.9957936666667"
$ws.Range("N20").Value = [double]"794.9873809999999"
$ws.Range("O20").Value = [double]"0.733559029746061"
$ws.Range("P20").Value = [double]"0.733559029746061"
$ws.Range("Q20").Value = [double]"73207.60307549436"
$ws.Range("R20").Value = [double]"658868.4276794492"
$ws.Range("S20").Value = [double]"0.4662371275940438"
$ws.Range("T20").Value = [double]"0.466237127594044"

$ws.Range("G21").Value = [double]"276.259491"
$ws.Range("H21").Value = [double]"828.7784730000001"
$ws.Range("I21").Value = [double]"0.6355822894790661"
$ws.Range("J21").Value = [double]"0.6355822894790663"
$ws.Range("M21").Value = [double]"2.197830333333333"
$ws.Range("N21").Value = [double]"6.593490999999999"
$ws.Range("O21").Value = [double]"0.006084014634943477"
$ws.Range("P21").Value = [double]"0.006084014634943477"
$ws.Range("Q21").Value = [double]"607.171489191027"
$ws.Range("R21").Value = [double]"5464.543402719243"
$ws.Range("S21").Value = [double]"0.00386689195090152"
$ws.Range("T21").Value = [double]"0.003866891950901521"

$ws.Range("G22").Value = [double]"71.006541"
$ws.Range("H22").Value = [double]"213.019623"
$ws.Range("I22").Value = [double]"0.1633627128371462"
$ws.Range("J22").Value = [double]"0.1633627128371462"
$ws.Range("K22").Value = [double]"2"
$ws.Range("L22").Value = [double]"0.6666666666666666"
$ws.Range("M22").Value = [double]"0.926994"
$ws.Range("N22").Value = [double]"2.780982"
$ws.Range("O22").Value = [double]"0.002566096653125693"
$ws.Range("P22").Value = [double]"0.002566096653125693"
$ws.Range("Q22").Value = [double]"65.82263746775399"
$ws.Range("R22").Value = [double]"592.403737209786"
$ws.Range("S22").Value = [double]"0.0004192045106569346"
$ws.Range("T22").Value = [double]"0.0004192045106569347"

$ws.Range("G23").Value = [double]"71.006541"
$ws.Range("H23").Value = [double]"213.019623"
$ws.Range("I23").Value = [double]"0.1633627128371462"
$ws.Range("J23").Value = [double]"0.1633627128371462"
$ws.Range("M23").Value = [double]"93.12610233333334"
$ws.Range("N23").Value = [double]"279.378307"
$ws.Range("O23").Value = [double]"0.2577908589658698"
$ws.Range("P23").Value = [double]"0.2577908589658698"
$ws.Range("Q23").Value = [double]"6612.562403502029"
$ws.Range("R23").Value = [double]"59513.06163151826"
$ws.Range("S23").Value = [double]"0.04211341406528264"
$ws.Range("T23").Value = [double]"0.04211341406528265"

$ws.Range("G24").Value = [double]"71.006541"
$ws.Range("H24").Value = [double]"213.019623"
$ws.Range("I24").Value = [double]"0.1633627128371462"
$ws.Range("J24").Value = [double]"0.1633627128371462"
$ws.Range("M24").Value = [double]"264.9957936666667"
$ws.Range("N24").Value = [double]"794.9873809999999"
$ws.Range("O24").Value = [double]"0.733559029746061"
$ws.Range("P24").Value = [double]"0.733559029746061"
$ws.Range("Q24").Value = [double]"18816.43468781971"
$ws.Range("R24").Value = [double]"169347.9121903773"
$ws.Range("S24").Value = [double]"0.1198361931255014"
$ws.Range("T24").Value = [double]"0.1198361931255014"

$ws.Range("G25").Value = [double]"71.006541"
$ws.Range("H25").Value = [double]"213.019623"
$ws.Range("I25").Value = [double]"0.1633627128371462"
$ws.Range("J25").Value = [double]"0.1633627128371462"
$ws.Range("M25").Value = [double]"2.197830333333333"
$ws.Range("N25").Value = [double]"6.593490999999999"
$ws.Range("O25").Value = [double]"0.006084014634943477"
$ws.Range("P25").Value = [double]"0.006084014634943477"
$ws.Range("Q25").Value = [double]"156.060329674877"
$ws.Range("R25").Value = [double]"1404.542967073893"
$ws.Range("S25").Value = [double]"0.0009939011357052662"
$ws.Range("T25").Value = [double]"0.0009939011357052662"
